$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JS-SPA-Self-Evaluation-Protocol")

# Days Commit in GitHub (C9): 53 -> 64
$ws.Range("C9").Value = 64

# Implement Filtering by Status score (C15): 4 -> 5
$ws.Range("C15").Value = 5

# Admin Home Screen / Admin Approve Ad / Admin Reject Ad now answered "Yes"
$ws.Range("C34").Value = "Yes"
$ws.Range("C35").Value = "Yes"
$ws.Range("C36").Value = "Yes"

# Leave the sheet's selection on C36, matching the saved view state
$ws.Range("C36").Select()
